# "add lock screen in step 1 for mitigate mutiples requests"
# Update the printed/header date on the first sheet and refresh the
# merged-cell regions (unmerge/re-merge) so the workbook's internal
# merge list reflects the latest edit order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date shown in A1 (2024-02-22 -> 2024-04-23)
$ws.Range("A1").Value = 45405

# Re-apply the merged cell ranges on the sheet (this refreshes their
# order in the saved workbook, matching a re-authored/regenerated sheet)
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:D1").Merge()

$ws.Range("A12:D12").UnMerge()
$ws.Range("A12:D12").Merge()

$ws.Range("B19:C19").UnMerge()
$ws.Range("B19:C19").Merge()

$ws.Range("A13:D13").UnMerge()
$ws.Range("A13:D13").Merge()

$ws.Range("A11:D11").UnMerge()
$ws.Range("A11:D11").Merge()

$ws.Range("B18:C18").UnMerge()
$ws.Range("B18:C18").Merge()

$ws.Range("B17:C17").UnMerge()
$ws.Range("B17:C17").Merge()

$ws.Range("A14:D14").UnMerge()
$ws.Range("A14:D14").Merge()

$wb.Save()
